$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4, shifting existing rows 4-12 down to 5-13.
$ws.Rows(4).Insert()

# Populate the new row 4 with the new data entry (weekly update).
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 44979
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100107
$ws.Cells.Item(4, 8).Value = "Otros"
$ws.Cells.Item(4, 9).Value = 100107011
$ws.Cells.Item(4, 10).Value = "Tuna"
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 250
$ws.Cells.Item(4, 14).Value = 29000
$ws.Cells.Item(4, 15).Value = 30000
$ws.Cells.Item(4, 16).Value = 29500
$ws.Cells.Item(4, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(4, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(4, 19).Value = 1475
$ws.Cells.Item(4, 20).Value = 20
